$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new PR log entry as row 30, matching the existing columns:
# ID | Title | Author | Approvers | Source / Target | Date
$ws.Range("A30").Value = 55
$ws.Range("B30").Value = "Edit1"
$ws.Range("C30").Value = "riya-morankar"
$ws.Range("D30").Value = "N/A"
$ws.Range("E30").Value = "edit1 to main"

# Force the Date column to remain plain text (matches the other rows,
# which store dates like "2025-06-17" as text, not as real dates).
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "2025-06-23"
